$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new column before column A, shifting the whole table from
#     A:G to B:H (values, number formats, column widths all move with it).
$ws.Columns.Item(1).Insert()

# --- Re-assert center/middle alignment across the (now shifted) table,
#     since this was already the look of the sheet before the edit.
$table = $ws.Range("B2:H7")
$table.HorizontalAlignment = -4108   # xlCenter
$table.VerticalAlignment = -4108     # xlCenter

# --- Shade the left-most (label) column of the table light gray, like a
#     row-header column.
$labelCol = $ws.Range("B2:B7")
$labelCol.Interior.Color = 14277081   # RGB(217,217,217) ~ theme bg1, -15% tint

# --- Thin grid lines across the whole table interior.
$table.Borders.Item(11).LineStyle = 1   # xlInsideVertical
$table.Borders.Item(11).Weight = 2      # xlThin
$table.Borders.Item(12).LineStyle = 1   # xlInsideHorizontal
$table.Borders.Item(12).Weight = 2      # xlThin

# --- Medium box border around the whole table.
$table.Borders.Item(7).LineStyle = 1    # xlEdgeLeft
$table.Borders.Item(7).Weight = -4138   # xlMedium
$table.Borders.Item(8).LineStyle = 1    # xlEdgeTop
$table.Borders.Item(8).Weight = -4138
$table.Borders.Item(9).LineStyle = 1    # xlEdgeBottom
$table.Borders.Item(9).Weight = -4138
$table.Borders.Item(10).LineStyle = 1   # xlEdgeRight
$table.Borders.Item(10).Weight = -4138

# --- Medium divider between the label column (B) and the rest of the table.
$labelCol.Borders.Item(10).LineStyle = 1  # xlEdgeRight
$labelCol.Borders.Item(10).Weight = -4138
$restCols = $ws.Range("C2:C7")
$restCols.Borders.Item(7).LineStyle = 1   # xlEdgeLeft
$restCols.Borders.Item(7).Weight = -4138

# --- Move the selection off the table, matching the saved cursor position.
$ws.Range("D25").Select()
